# Apply the commit's edit to distribution_circuits.xlsx
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1 edits ---
# (plotTts / Sheet2's header cells in the diff only shift shared-string
# indices because three now-unused strings are dropped below - the text
# itself is unchanged, so no edit is needed there.)

# Row 14: I14 8 -> 9
$ws1.Range("I14").Value = 9

# Row 14: N14, O14, P14 get values (X, O, X) - previously blank
$ws1.Range("N14").Value = "X"
$ws1.Range("O14").Value = "O"
$ws1.Range("P14").Value = "X"

# Row 14: Q14, R14 cleared (were "HARD"/1 leftover TODO markers)
$ws1.Range("Q14").ClearContents()
$ws1.Range("R14").ClearContents()

# Row 2: Q2, R2 cleared (were "DIFF"/"TODO ORDER" header leftovers)
$ws1.Range("Q2").ClearContents()
$ws1.Range("R2").ClearContents()

# Update the active selection on Sheet1 to L17 (was R15)
$ws1.Activate() | Out-Null
$ws1.Range("L17").Select() | Out-Null
